$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" sheet: insert a new row 2 for "2022-Q4", pushing the rest down.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)
$summary.Rows.Item(2).Insert()

# Give the new A2 the same style as the (shifted) row below it, then set value.
$summary.Range("A3").Copy($summary.Range("A2"))
$summary.Range("A2").Value = 0

# Clear the placeholder style picked up from Insert() on B2:D2 before filling.
$summary.Range("B2:D2").ClearFormats()
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 17
$summary.Range("D2").Value = 4.91

# ---------------------------------------------------------------------------
# 2) New "2022-Q4" worksheet: duplicate "2022-Q3" (keeps styles/columns),
#    rename, then overwrite its data with the 2022-Q4 numbers.
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($null, $summary)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# The template only has 9 data rows (2..10); extend style down to row 18.
$q4.Range("A10").Copy($q4.Range("A11:A18"))

# Build the 17x8 data array for the new 2022-Q4 sheet (rows 2..18, cols A..H)
$data = New-Object 'object[,]' 17,8
$data[0,0] = 0
$data[0,1] = "'" + '010041'
$data[0,2] = '嘉实港股优势混合A'
$data[0,3] = "'" + '45.52'
$data[0,4] = "'" + '93.43'
$data[0,5] = "'" + '3.37'
$data[0,6] = "'" + '1.5340'
$data[0,7] = 8
$data[1,0] = 1
$data[1,1] = "'" + '001878'
$data[1,2] = '嘉实沪港深精选股票'
$data[1,3] = "'" + '21.50'
$data[1,4] = "'" + '90.83'
$data[1,5] = "'" + '3.36'
$data[1,6] = "'" + '0.7224'
$data[1,7] = 10
$data[2,0] = 2
$data[2,1] = "'" + '010591'
$data[2,2] = '富国中国中小盘混合（QDII）美元'
$data[2,3] = "'" + '32.21'
$data[2,4] = "'" + '87.21'
$data[2,5] = "'" + '2.15'
$data[2,6] = "'" + '0.6925'
$data[2,7] = 9
$data[3,0] = 3
$data[3,1] = "'" + '100061'
$data[3,2] = '富国中国中小盘混合（QDII）人民币'
$data[3,3] = "'" + '32.21'
$data[3,4] = "'" + '87.21'
$data[3,5] = "'" + '2.15'
$data[3,6] = "'" + '0.6925'
$data[3,7] = 9
$data[4,0] = 4
$data[4,1] = "'" + '007291'
$data[4,2] = '汇丰晋信港股通双核策略混合'
$data[4,3] = "'" + '7.74'
$data[4,4] = "'" + '90.21'
$data[4,5] = "'" + '5.15'
$data[4,6] = "'" + '0.3986'
$data[4,7] = 4
$data[5,0] = 5
$data[5,1] = "'" + '002332'
$data[5,2] = '汇丰晋信沪港深股票A'
$data[5,3] = "'" + '4.46'
$data[5,4] = "'" + '90.48'
$data[5,5] = "'" + '4.52'
$data[5,6] = "'" + '0.2016'
$data[5,7] = 4
$data[6,0] = 6
$data[6,1] = "'" + '010042'
$data[6,2] = '嘉实港股优势混合C'
$data[6,3] = "'" + '5.86'
$data[6,4] = "'" + '93.43'
$data[6,5] = "'" + '3.37'
$data[6,6] = "'" + '0.1975'
$data[6,7] = 8
$data[7,0] = 7
$data[7,1] = "'" + '007182'
$data[7,2] = '万家沪港深蓝筹混合A'
$data[7,3] = "'" + '3.44'
$data[7,4] = "'" + '92.13'
$data[7,5] = "'" + '3.63'
$data[7,6] = "'" + '0.1249'
$data[7,7] = 8
$data[8,0] = 8
$data[8,1] = "'" + '002333'
$data[8,2] = '汇丰晋信沪港深股票C'
$data[8,3] = "'" + '1.79'
$data[8,4] = "'" + '90.48'
$data[8,5] = "'" + '4.52'
$data[8,6] = "'" + '0.0809'
$data[8,7] = 4
$data[9,0] = 9
$data[9,1] = "'" + '013009'
$data[9,2] = '万家港股通精选混合A'
$data[9,3] = "'" + '2.56'
$data[9,4] = "'" + '84.90'
$data[9,5] = "'" + '3.06'
$data[9,6] = "'" + '0.0783'
$data[9,7] = 10
$data[10,0] = 10
$data[10,1] = "'" + '003993'
$data[10,2] = '前海开源沪港深核心驱动灵活配置混合'
$data[10,3] = "'" + '1.16'
$data[10,4] = "'" + '72.10'
$data[10,5] = "'" + '5.13'
$data[10,6] = "'" + '0.0595'
$data[10,7] = 8
$data[11,0] = 11
$data[11,1] = "'" + '006537'
$data[11,2] = '恒生前海港股通精选混合'
$data[11,3] = "'" + '1.10'
$data[11,4] = "'" + '91.14'
$data[11,5] = "'" + '3.66'
$data[11,6] = "'" + '0.0403'
$data[11,7] = 6
$data[12,0] = 12
$data[12,1] = "'" + '007183'
$data[12,2] = '万家沪港深蓝筹混合C'
$data[12,3] = "'" + '0.86'
$data[12,4] = "'" + '92.13'
$data[12,5] = "'" + '3.63'
$data[12,6] = "'" + '0.0312'
$data[12,7] = 8
$data[13,0] = 13
$data[13,1] = "'" + '006477'
$data[13,2] = '中邮沪港深精选混合'
$data[13,3] = "'" + '0.67'
$data[13,4] = "'" + '94.24'
$data[13,5] = "'" + '4.59'
$data[13,6] = "'" + '0.0308'
$data[13,7] = 10
$data[14,0] = 14
$data[14,1] = "'" + '013010'
$data[14,2] = '万家港股通精选混合C'
$data[14,3] = "'" + '0.77'
$data[14,4] = "'" + '84.90'
$data[14,5] = "'" + '3.06'
$data[14,6] = "'" + '0.0236'
$data[14,7] = 10
$data[15,0] = 15
$data[15,1] = "'" + '013182'
$data[15,2] = '安信港股通精选混合C'
$data[15,3] = "'" + '0.12'
$data[15,4] = "'" + '69.28'
$data[15,5] = "'" + '3.12'
$data[15,6] = "'" + '0.0037'
$data[15,7] = 7
$data[16,0] = 16
$data[16,1] = "'" + '013181'
$data[16,2] = '安信港股通精选混合A'
$data[16,3] = "'" + '0.02'
$data[16,4] = "'" + '69.28'
$data[16,5] = "'" + '3.12'
$data[16,6] = "'" + '0.0006'
$data[16,7] = 7

# Write the new data as a single block assignment (avoids per-cell overhead).
$q4.Range("A2:H18").Value = $data

# Strings that must stay text (fund code / decimal-looking numbers as text)
# picked up the "Text" number format style from the auto-text coercion above;
# reset those ranges to the sheet's normal (unstyled) look, matching the
# other quarters where only the header row and column A carry a style.
$q4.Range("B2:B18").Style = "Normal"
$q4.Range("D2:G18").Style = "Normal"

Write-Host "2022-Q4 sheet populated"
